# Fruta / hortaliza, semanal
# Insert two new weekly records for "Agrícola del Norte S.A. de Arica - Pera"
# as rows 15 and 16, pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 15 (existing rows 15.. shift to 17..)
$ws.Rows("15:16").Insert()

# ---- New row 15 ----
$r = 15
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($r, 4).Value = 45177
$ws.Cells.Item($r, 5).Value = 15
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100104
$ws.Cells.Item($r, 8).Value = "Frutos de pepita"
$ws.Cells.Item($r, 9).Value = 100104005
$ws.Cells.Item($r, 10).Value = "Pera"
$ws.Cells.Item($r, 11).Value = "Packham's Triumph"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 420
$ws.Cells.Item($r, 14).Value = 22000
$ws.Cells.Item($r, 15).Value = 23000
$ws.Cells.Item($r, 16).Value = 22286
$ws.Cells.Item($r, 17).Value = "`$/caja 20 kilos granel"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 1114
$ws.Cells.Item($r, 20).Value = 20

# ---- New row 16 ----
$r = 16
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($r, 4).Value = 45177
$ws.Cells.Item($r, 5).Value = 15
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100104
$ws.Cells.Item($r, 8).Value = "Frutos de pepita"
$ws.Cells.Item($r, 9).Value = 100104005
$ws.Cells.Item($r, 10).Value = "Pera"
$ws.Cells.Item($r, 11).Value = "Packham's Triumph"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 200
$ws.Cells.Item($r, 14).Value = 20000
$ws.Cells.Item($r, 15).Value = 21000
$ws.Cells.Item($r, 16).Value = 20500
$ws.Cells.Item($r, 17).Value = "`$/caja 20 kilos granel"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 1025
$ws.Cells.Item($r, 20).Value = 20
